$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11 (Rating, Precision, Recall, F1-Score, Support)
$data = @(
    @("A",   "0.9242", "0.9112", "0.9176", "214"),
    @("AA",  "0.8947", "0.8095", "0.8500", "42"),
    @("AAA", "0.9231", "0.8889", "0.9057", "27"),
    @("B",   "0.9362", "0.9565", "0.9462", "138"),
    @("BB",  "0.9701", "0.9799", "0.9750", "298"),
    @("BBB", "0.9548", "0.9635", "0.9592", "329"),
    @("C",   "1.0000", "1.0000", "1.0000", "7"),
    @("CC",  "0.0000", "0.0000", "0.0000", "0"),
    @("CCC", "0.9688", "0.9118", "0.9394", "34"),
    @("D",   "1.0000", "1.0000", "1.0000", "3")
)

$rowIndex = 2
foreach ($rowData in $data) {
    $ws.Range("A$rowIndex").Value = "'" + $rowData[0]
    $ws.Range("B$rowIndex").Value = "'" + $rowData[1]
    $ws.Range("C$rowIndex").Value = "'" + $rowData[2]
    $ws.Range("D$rowIndex").Value = "'" + $rowData[3]
    $ws.Range("E$rowIndex").Value = "'" + $rowData[4]

    $rowIndex++
}
